# "update der kampfhandlung methode"
#
# - rename the stat label "LeP" (life points) in A16 to "LP"
# - clear the (unused, identical) extra cell style that was applied to
#   the whole sheet, so cells go back to the default/general style
# - move the active cell / selection to D20

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "LeP" label to "LP"
$ws.Range("A16").Value2 = "LP"

# Drop the extra cell formatting (s="1") that was applied to every
# used cell; apply to the two contiguous blocks of actual data so we
# don't touch the empty row 15 in between (which would otherwise get
# materialized into sheetData).
$ws.Range("A1:B14").Style = "Standard"
$ws.Range("A16:B22").Style = "Standard"

# Update the selected/active cell shown in the sheet view
$ws.Range("D20").Select()
